$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C9").Value = "ip:"
$ws.Range("D9").Value = "10.10.69.69"

$ws.Range("D11").Select()
